# Update ObjTables header row (A1) on every worksheet: the legacy
# "id='ClassName'" attribute has been renamed to "class='ClassName'"
# by the obj_tables library. Apply this rename wherever it appears.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A1")
    $text = $cell.Value2

    if ($text -ne $null -and $text -like "*!!ObjTables*" -and $text -like "*id='*") {
        $newText = $text -replace "id='", "class='"
        $cell.Value = $newText
    }
}
